# Fruta / hortaliza, semanal
# Update rows 3,5,7,8,9,10,12,13 on the active sheet with new Fecha (D),
# Volumen (M), Precio minimo/maximo/promedio (N/O/P), Origen (R) and
# Precio $/Kg (S) values, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value2 = 44999   # D3 Fecha
$ws.Cells.Item(3, 13).Value2 = 25     # M3 Volumen

# Row 5
$ws.Cells.Item(5, 4).Value2 = 44215   # D5 Fecha
$ws.Cells.Item(5, 13).Value2 = 65     # M5 Volumen
$ws.Cells.Item(5, 14).Value2 = 2800   # N5 Precio minimo
$ws.Cells.Item(5, 15).Value2 = 2800   # O5 Precio maximo
$ws.Cells.Item(5, 16).Value2 = 2800   # P5 Precio promedio ponderado
$ws.Cells.Item(5, 19).Value2 = 2800   # S5 Precio $/Kg

# Row 7
$ws.Cells.Item(7, 4).Value2 = 44998   # D7 Fecha
$ws.Cells.Item(7, 13).Value2 = 20     # M7 Volumen
$ws.Cells.Item(7, 14).Value2 = 2500   # N7 Precio minimo
$ws.Cells.Item(7, 15).Value2 = 2500   # O7 Precio maximo
$ws.Cells.Item(7, 16).Value2 = 2500   # P7 Precio promedio ponderado
$ws.Cells.Item(7, 19).Value2 = 2500   # S7 Precio $/Kg

# Row 8
$ws.Cells.Item(8, 4).Value2 = 44175   # D8 Fecha
$ws.Cells.Item(8, 13).Value2 = 40     # M8 Volumen
$ws.Cells.Item(8, 14).Value2 = 5000   # N8 Precio minimo
$ws.Cells.Item(8, 15).Value2 = 5000   # O8 Precio maximo
$ws.Cells.Item(8, 16).Value2 = 5000   # P8 Precio promedio ponderado
$ws.Cells.Item(8, 18).Value = "Provincia de Curicó"  # R8 Origen
$ws.Cells.Item(8, 19).Value2 = 5000   # S8 Precio $/Kg

# Row 9
$ws.Cells.Item(9, 4).Value2 = 44176   # D9 Fecha
$ws.Cells.Item(9, 13).Value2 = 20     # M9 Volumen
$ws.Cells.Item(9, 14).Value2 = 3000   # N9 Precio minimo
$ws.Cells.Item(9, 15).Value2 = 3000   # O9 Precio maximo
$ws.Cells.Item(9, 16).Value2 = 3000   # P9 Precio promedio ponderado
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"  # R9 Origen
$ws.Cells.Item(9, 19).Value2 = 3000   # S9 Precio $/Kg

# Row 10
$ws.Cells.Item(10, 4).Value2 = 44551   # D10 Fecha
$ws.Cells.Item(10, 13).Value2 = 120    # M10 Volumen
$ws.Cells.Item(10, 14).Value2 = 4500   # N10 Precio minimo
$ws.Cells.Item(10, 15).Value2 = 4500   # O10 Precio maximo
$ws.Cells.Item(10, 16).Value2 = 4500   # P10 Precio promedio ponderado
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"  # R10 Origen
$ws.Cells.Item(10, 19).Value2 = 4500   # S10 Precio $/Kg

# Row 12
$ws.Cells.Item(12, 4).Value2 = 44323   # D12 Fecha
$ws.Cells.Item(12, 14).Value2 = 3200   # N12 Precio minimo
$ws.Cells.Item(12, 15).Value2 = 3200   # O12 Precio maximo
$ws.Cells.Item(12, 16).Value2 = 3200   # P12 Precio promedio ponderado
$ws.Cells.Item(12, 18).Value = "Región de La Araucanía"  # R12 Origen
$ws.Cells.Item(12, 19).Value2 = 3200   # S12 Precio $/Kg

# Row 13
$ws.Cells.Item(13, 4).Value2 = 44592   # D13 Fecha
$ws.Cells.Item(13, 13).Value2 = 5      # M13 Volumen
$ws.Cells.Item(13, 14).Value2 = 7500   # N13 Precio minimo
$ws.Cells.Item(13, 15).Value2 = 7500   # O13 Precio maximo
$ws.Cells.Item(13, 16).Value2 = 7500   # P13 Precio promedio ponderado
$ws.Cells.Item(13, 18).Value = "Región de La Araucanía"  # R13 Origen
$ws.Cells.Item(13, 19).Value2 = 7500   # S13 Precio $/Kg
